$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Lot0): Status OCUPAT -> LIBER, Date cleared to 0, BusySpaces 3 -> 2, FreeSpaces 22 -> 23
$ws.Range("B2").Value = "LIBER"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 23

# Row 3 (Lot1): Date updated, BusySpaces 3 -> 2, FreeSpaces 22 -> 23
$ws.Range("C3").Value = "04-24-2021 16:56:56"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 23

# Row 4 (Lot2): Date updated, BusySpaces 3 -> 2, FreeSpaces 22 -> 23
$ws.Range("C4").Value = "04-24-2021 16:56:56"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 23
